$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Configs" right before "Possible policies combinations"
$target = $wb.Worksheets.Item("Possible policies combinations")
$new = $wb.Worksheets.Add($target)
$new.Name = "Configs"

# Populate the new sheet with the config key/value pairs
$new.Range("A1").Value = "population"
$new.Range("B1").Value = 4200

$new.Range("A2").Value = "vegetable_demand"
$new.Range("B2").Value = "0.125*365"

$new.Range("A3").Value = "agri_land_availability"
$new.Range("B3").Value = 27000

$new.Range("A4").Value = "energy_land_availabilty"
$new.Range("B4").Value = 140000

$new.Range("A5").Value = "electricity_demand"
$new.Range("B5").Value = 43200000

$new.Range("A6").Value = "error_buffer"
$new.Range("B6").Value = 0

$new.Range("A7").Value = "grid_lenght"
$new.Range("B7").Value = 40

# Column A autofit width (matches bestFit width seen in target sheet)
$new.Columns.Item(1).AutoFit() | Out-Null

# Select C5 on the new sheet, matching target selection
$new.Range("C5").Select() | Out-Null

# Activate the new sheet (it becomes the selected tab)
$new.Activate() | Out-Null
